$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @{
    2  = @(2.828655049581922, 1.773799324762276, 1.493186813606976, 1.320195973509991)
    3  = @(2.883834815612763, 1.790830877570518, 1.500165029430642, 1.325315907228027)
    4  = @(3.01732006345383,  1.860444782376231, 1.549274810175597, 1.36634262392983)
    5  = @(3.380719459040112, 1.998626836058858, 1.601553576152952, 1.409259220543114)
    6  = @(3.364883251118721, 1.996492174193518, 1.602459148828121, 1.409908246290728)
    7  = @(3.655791957161804, 2.082638963899223, 1.641837099424327, 1.414435199910633)
    8  = @(1.61754885472335,  1.331268593591207, 1.25963316010294,  1.286428854980599)
    9  = @(3.03186021144269,  1.862728241502911, 1.549120835106842, 1.366023141857575)
    10 = @(3.678504460331955, 2.084409928598445, 1.639123914208648, 1.409810408930862)
    11 = @(3.781325438004655, 2.106358068430188, 1.646314611292481, 1.414360880002128)
    12 = @(3.573545316688634, 2.058589535148215, 1.6308334743146,   1.405296775565357)
    13 = @(3.657686037493022, 2.079739547593903, 1.636853583344017, 1.408286739559427)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
}
